$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header for column B
$ws.Range("B1").Value = "CombinedVaderSentiment"

# Update sentiment values in column B (rows 2-40), leaving unchanged rows untouched
$ws.Range("B2").Value = 0.959
$ws.Range("B3").Value = 0.56335
$ws.Range("B4").Value = 0.5781000000000001
$ws.Range("B5").Value = 0.852
$ws.Range("B6").Value = 0.995
$ws.Range("B7").Value = 0.9924999999999999
$ws.Range("B8").Value = 0.995
$ws.Range("B9").Value = 0.995
$ws.Range("B10").Value = 0.986
$ws.Range("B11").Value = 0.8704499999999999
$ws.Range("B12").Value = 0.9975000000000001
$ws.Range("B13").Value = 1.1598
$ws.Range("B14").Value = 0.998
$ws.Range("B15").Value = 0.19515
$ws.Range("B16").Value = 0.86225
$ws.Range("B17").Value = 0.957
$ws.Range("B18").Value = 0.998
$ws.Range("B20").Value = 0.9804999999999999
$ws.Range("B21").Value = 0.9895
$ws.Range("B22").Value = 1.0248
$ws.Range("B25").Value = 0.86465
$ws.Range("B26").Value = 0.996
$ws.Range("B28").Value = 0.957
$ws.Range("B29").Value = 0.988
$ws.Range("B30").Value = 0.993
$ws.Range("B31").Value = 0.98
$ws.Range("B32").Value = 0.998
$ws.Range("B34").Value = 0.995
$ws.Range("B35").Value = 0.90655
$ws.Range("B36").Value = 0.7800499999999999
$ws.Range("B37").Value = 0.9955000000000001
$ws.Range("B38").Value = 0.99
$ws.Range("B39").Value = 0.9924999999999999
